{"js": "// Analysis of the provided diff: every hunk only reorders XML attributes\n// (and namespace declarations) alphabetically, normalizes whitespace\n// inside a couple of base64-encoded VML `o:gfxdata` blobs, and drops\n// transient `w:rsid*` bookkeeping attributes that Word regenerates on\n// every save. None of the w:val/r:id/measurement values, text runs,\n// paragraph/style structure, relationships, or any other semantic\n// content differ between the \"before\" and \"after\" OOXML. In other\n// words this commit (\"Moving from 2.0.2 to 2.0.3\") is purely a\n// re-serialization done by the authoring tool's updated XML writer -\n// there is no actual document edit to reproduce through the Word API.\n//\n// Touch the document read-only (load the body, sync) so the script\n// demonstrably exercises the API without mutating any content.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# Analysis of the provided diff: every hunk only reorders XML attributes\n# (and namespace declarations) alphabetically, normalizes whitespace\n# inside a couple of base64-encoded VML o:gfxdata blobs, and drops\n# transient w:rsid* bookkeeping attributes that Word regenerates on\n# every save. None of the w:val/r:id/measurement values, text runs,\n# paragraph/style structure, relationships, or any other semantic\n# content differ between the \"before\" and \"after\" OOXML. In other\n# words this commit (\"Moving from 2.0.2 to 2.0.3\") is purely a\n# re-serialization done by the authoring tool's updated XML writer -\n# there is no actual document edit to reproduce through the Word\n# object model.\n#\n# Touch the document read-only (read the body text) so the script\n# demonstrably exercises the object model without mutating any content.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
